# Apply the "Added all current Data" edit:
#  - Insert a new column B ("segments") containing the original segment
#    labels (which used to live in column A).
#  - Column A becomes a plain numeric index (0-based row order).
#  - The old PercActivations / PercSegmentAreas columns shift from B/C to C/D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new (empty) column before the existing column B.
#    This shifts the old B (PercActivations) -> C and old C (PercSegmentAreas) -> D,
#    carrying their values/formats/header style along automatically.
$ws.Columns("B").Insert()

# 2. New column B header - give it the same bold/bordered/centered header
#    style as the neighbouring header cells, then set its text.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1").Value2 = "segments"

# 3. Move the segment-name labels that are still sitting in column A
#    (rows 2-20) over to the new column B, then replace column A with a
#    simple 0-based numeric index.
$lastRow = 20
for ($r = 2; $r -le $lastRow; $r++) {
    $label = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 2).Value2 = $label
    $ws.Cells.Item($r, 1).Value2 = $r - 2
}

# 4. The insert copied column A's header-row style onto the new column B
#    for the data rows too; clear that back to the default (unstyled) look
#    used by the original PercActivations/PercSegmentAreas data columns.
$ws.Range("B2:B20").ClearFormats()
